# Update the base model data: the "Plain English" column header was
# stored as a stray/duplicate shared string; retype it (with an
# underscore, "Plain_English") in every lookup sheet that uses it so
# Excel folds it back into its natural place in the shared-strings table.
$wb = $excel.ActiveWorkbook

$wb.Sheets.Item("Fuel_to_Code").Range("C1").Value = "Plain_English"
$wb.Sheets.Item("VehFuel_to_Code").Range("C1").Value = "Plain_English"
$wb.Sheets.Item("Tech_to_Code").Range("C1").Value = "Plain_English"
$wb.Sheets.Item("Dem_to_Code").Range("B1").Value = "Plain_English"

# Also clean the folder: leave the selections/active sheet the way the
# author's session ended up (Dem_to_Code active, with the other lookup
# sheets' cursors parked where they were left).
$wb.Sheets.Item("Fuel_to_Code").Range("D8").Select()
$wb.Sheets.Item("VehFuel_to_Code").Range("A10:XFD11").Select()
$wb.Sheets.Item("Tech_to_Code").Range("A1:D13").Select()

$ws7 = $wb.Sheets.Item("Dem_to_Code")
$ws7.Activate()
$ws7.Range("F13").Select()
